$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns (G:J) for the new National/Senate percentage columns.
# This shifts the existing G:AE block (district-level results) to K:AI,
# and Excel automatically rewrites the formulas that referenced the old
# column letters (e.g. H1 -> L1).
$ws.Range("G1:J1").EntireColumn.Insert()

# Header row (row 2): new column headers for the inserted columns.
# Entered in this order so the shared-string table ends up with the same
# ids as the source workbook (I2/J2 created first, then G2/H2).
$ws.Range("I2").Value = "% DEM Sen"
$ws.Range("J2").Value = "% REP Sen"
$ws.Range("G2").Value = "% DEM Nat"
$ws.Range("H2").Value = "% REP Nat"

# Row 3 (Idaho) - national/senate percentage data for the new columns.
$ws.Range("G3").Value = 0.11
$ws.Range("H3").Value = 0.1
$ws.Range("I3").Value = 0.1
$ws.Range("J3").Value = 0.11

# Row 4 (Nevada) - national/senate percentage data for the new columns.
$ws.Range("G4").Value = 0.1
$ws.Range("H4").Value = 0.11
$ws.Range("I4").Value = 0.11
$ws.Range("J4").Value = 0.1

# A bolded (but otherwise empty) band further down the sheet.
$ws.Range("G28:M28").Font.Bold = $true

# Leave the same cell selected as in the edited workbook.
[void]$ws.Range("J11").Select()
